$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Region_selection")

# Insert a new row at row 2 (pushes existing data down by one row)
$ws.Rows.Item(2).Insert()

# Fill in the new "World" entry
$ws.Range("A2").Value = "World"
$ws.Range("B2").Value = 1

# Make the new row's A/B cells explicitly non-bold (matches the rest of the
# data rows) while C2 picks up the bold style used by the header row's C1
$ws.Range("A2:B2").Font.Bold = $false
$ws.Range("C2").Font.Bold = $true

$ws.Range("A2").Select()
